$d = $word.ActiveDocument

# 1. Wrap the whole document body in a "_GoBack" bookmark (id 0), matching
#    the range Word marks as "last edit location" after a round of edits.
$full = $d.Content
$full.End = $full.End - 1
$d.Bookmarks.Add("_GoBack", $full)

# 2. Recognize "Heading1" as true outline-level-1 heading (outlineLvl 0)
#    instead of outline level 2 (outlineLvl 1).
$hs = $d.Styles("Heading1")
$hs.ParagraphFormat.OutlineLevel = 1

# 3. Normalize section page setup: explicit header/footer distance and
#    gutter, plus default column spacing -- matches Word re-serializing the
#    sectPr with its usual defaults spelled out.
$ps = $d.Sections(1).PageSetup
$ps.HeaderDistance = 36
$ps.FooterDistance = 36
$ps.Gutter = 0
$ps.TextColumns.Spacing = 36
